$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Visual Upgrades" detail (row 13, col D): the old text had two bullet
#    lines — the first (running animation) becomes its own new row below
#    ("Player Animations"), so this cell keeps only the second bullet line.
# ---------------------------------------------------------------------------
$ws.Range("D13").Value = "(-> improved background art and ui button art for the gameplay buttons"

# ---------------------------------------------------------------------------
# 2) Split the big "AI Attacking" rich-text cell (row 10, col D): the trailing
#    paragraph about the AI coach becomes its own (green) run instead of
#    sharing the default/black run with the blank line above it.
# ---------------------------------------------------------------------------
$cell = $ws.Range("D10")
$full = $cell.Characters(1, 4000).Text
$marker = "(-> have the AI coach"
$idx0 = $full.IndexOf($marker)
if ($idx0 -ge 0) {
    $startChar = $idx0          # 1-based char index of the newline right before the marker
    $length = $full.Length - $idx0 + 1
    $chars = $cell.Characters($startChar, $length)
    $chars.Font.Color = 5296274   # RGB(0x92,0xD0,0x50) -> green, matches the other bullet runs
}

# ---------------------------------------------------------------------------
# 3) Insert a new row 14 ("Player Animations") between the existing "Visual
#    Upgrades" row (13) and the blank spacer row that precedes "AI will
#    react to player pawns positions" (old row 15, now row 16).
# ---------------------------------------------------------------------------
$ws.Rows("14").Insert()
$ws.Range("B14").Value = "Player Animations"
$ws.Range("C14").Value = "Garrett"
$ws.Range("D14").Value = "(-> add animations for both player and AI`n(-> hook up animations to each specific phase"
$ws.Range("D14").WrapText = $true

# The blank spacer row that used to sit right under the big "AI will react.."
# cell (old row 16) is dropped entirely -- after the insert above it has been
# pushed down to row 17, so remove it now, which re-aligns every row below
# back to its original row number.
$ws.Rows("17").Delete()

# ---------------------------------------------------------------------------
# 4) Update the view state to match where the author left the selection.
# ---------------------------------------------------------------------------
$ws.Range("D16").Select()
$excel.ActiveWindow.ScrollRow = 11
